$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.316689729690552
$ws.Range("B1").Value = 1.865329742431641
$ws.Range("C1").Value = 2.692427158355713
$ws.Range("D1").Value = 4.921469688415527
$ws.Range("E1").Value = 1.144118905067444
